$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 386237.66
$ws.Range("J17").Value = 386237.66
$ws.Range("L17").Value = 1158712.98
$ws.Range("N17").Value = -1159048.98
$ws.Range("H18").Value = 13458.4
$ws.Range("I18").Value = 17500.285
$ws.Range("K18").Value = 17500.285
$ws.Range("M18").Value = -17216.285
$ws.Range("H98").Value = 2746.4285
$ws.Range("I98").Value = 1405.5454
$ws.Range("J98").Value = 7663
$ws.Range("K98").Value = 1405.5454
$ws.Range("L98").Value = 7663
$ws.Range("M98").Value = 92.45460000000003
$ws.Range("N98").Value = -10659
$ws.Range("H111").Value = 0
$ws.Range("I111").Value = 0
$ws.Range("K111").Value = 0
$ws.Range("M111").Value = ""
$ws.Range("H112").Value = 2014.9333
$ws.Range("J112").Value = 2196.5
$ws.Range("L112").Value = 6589.5
$ws.Range("N112").Value = -8805.5
$ws.Range("H122").Value = 2746.4285
$ws.Range("I122").Value = 1405.5454
$ws.Range("J122").Value = 7663
$ws.Range("K122").Value = 4216.6362
$ws.Range("L122").Value = 22989
$ws.Range("M122").Value = -1766.6362
$ws.Range("N122").Value = -27889
$ws.Range("H138").Value = 1765.3871
$ws.Range("I138").Value = 1239.9445
$ws.Range("J138").Value = 2492.923
$ws.Range("K138").Value = 3719.8335
$ws.Range("L138").Value = 7478.768999999999
$ws.Range("M138").Value = 1420.1665
$ws.Range("N138").Value = -17758.769

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 7456.0454
$ws.Range("I45").Value = 10553.417
$ws.Range("K45").Value = 10553.417
$ws.Range("M45").Value = -10176.417
$ws.Range("H57").Value = 4942.6665
$ws.Range("I57").Value = 4942.6665
$ws.Range("K57").Value = 4942.6665
$ws.Range("M57").Value = -4458.6665
$ws.Range("H61").Value = 1813.091
$ws.Range("I61").Value = 1694.4
$ws.Range("K61").Value = 1694.4
$ws.Range("M61").Value = -1482.4
$ws.Range("H136").Value = 1813.091
$ws.Range("I136").Value = 1694.4
$ws.Range("K136").Value = 5083.200000000001
$ws.Range("M136").Value = -2533.200000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 8414.5
$ws.Range("I20").Value = 8509.2
$ws.Range("K20").Value = 8509.2
$ws.Range("M20").Value = -8262.2
$ws.Range("H76").Value = 23521.334
$ws.Range("J76").Value = 23782
$ws.Range("L76").Value = 23782
$ws.Range("N76").Value = -24412
$ws.Range("H79").Value = 23521.334
$ws.Range("J79").Value = 23782
$ws.Range("L79").Value = 23782
$ws.Range("N79").Value = -25966
$ws.Range("H122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("N122").Value = ""

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2189.625
$ws.Range("I132").Value = 2133.12
$ws.Range("J132").Value = 2391.4285
$ws.Range("K132").Value = 6399.36
$ws.Range("L132").Value = 7174.2855
$ws.Range("M132").Value = -3869.36
$ws.Range("N132").Value = -12234.2855
$ws.Range("H134").Value = 3920.9092
$ws.Range("I134").Value = 3313
$ws.Range("K134").Value = 9939
$ws.Range("M134").Value = -7404

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 8916.277
$ws.Range("I56").Value = 8916.277
$ws.Range("K56").Value = 8916.277
$ws.Range("M56").Value = -8386.277
$ws.Range("H113").Value = 901.3333
$ws.Range("J113").Value = 902
$ws.Range("L113").Value = 2706
$ws.Range("N113").Value = -7046

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 35216
$ws.Range("J52").Value = 35216
$ws.Range("L52").Value = 35216
$ws.Range("N52").Value = -35734

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 29466.166
$ws.Range("I7").Value = 41337.125
$ws.Range("K7").Value = 41337.125
$ws.Range("M7").Value = -41225.125
$ws.Range("H18").Value = 8000
$ws.Range("I18").Value = 0
$ws.Range("J18").Value = 8000
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = 8000
$ws.Range("M18").Value = ""
$ws.Range("N18").Value = -8344
$ws.Range("H42").Value = 22998.8
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 22998.8
$ws.Range("K42").Value = 0
$ws.Range("L42").Value = 22998.8
$ws.Range("M42").Value = ""
$ws.Range("N42").Value = -24124.8
$ws.Range("H46").Value = 28359.938
$ws.Range("I46").Value = 71910.336
$ws.Range("K46").Value = 71910.336
$ws.Range("M46").Value = -71722.336
$ws.Range("H49").Value = 22998.8
$ws.Range("I49").Value = 0
$ws.Range("J49").Value = 22998.8
$ws.Range("K49").Value = 0
$ws.Range("L49").Value = 22998.8
$ws.Range("M49").Value = ""
$ws.Range("N49").Value = -23292.8
$ws.Range("H61").Value = 15879.125
$ws.Range("I61").Value = 1756.1578
$ws.Range("J61").Value = 69546.4
$ws.Range("K61").Value = 1756.1578
$ws.Range("L61").Value = 69546.4
$ws.Range("M61").Value = -1554.1578
$ws.Range("N61").Value = -69950.4
$ws.Range("H113").Value = 15879.125
$ws.Range("I113").Value = 1756.1578
$ws.Range("J113").Value = 69546.4
$ws.Range("K113").Value = 1756.1578
$ws.Range("L113").Value = 69546.4
$ws.Range("M113").Value = 413.8422
$ws.Range("N113").Value = -73886.4
$ws.Range("H122").Value = 137603.06
$ws.Range("I122").Value = 289008
$ws.Range("K122").Value = 867024
$ws.Range("M122").Value = -864574
$ws.Range("H126").Value = 29466.166
$ws.Range("I126").Value = 41337.125
$ws.Range("K126").Value = 124011.375
$ws.Range("M126").Value = -121541.375
$ws.Range("H133").Value = 119998
$ws.Range("J133").Value = 119998
$ws.Range("L133").Value = 119998
$ws.Range("N133").Value = -125058

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 3522
$ws.Range("I81").Value = 2033.625
$ws.Range("J81").Value = 6498.75
$ws.Range("K81").Value = 4067.25
$ws.Range("L81").Value = 12997.5
$ws.Range("M81").Value = -3006.25
$ws.Range("N81").Value = -15119.5
$ws.Range("H84").Value = 3522
$ws.Range("I84").Value = 2033.625
$ws.Range("J84").Value = 6498.75
$ws.Range("K84").Value = 20336.25
$ws.Range("L84").Value = 64987.5
$ws.Range("M84").Value = -15032.25
$ws.Range("N84").Value = -75595.5
$ws.Range("H88").Value = 37585.5
$ws.Range("J88").Value = 25000
$ws.Range("L88").Value = 25000
$ws.Range("N88").Value = -25812
$ws.Range("H91").Value = 37585.5
$ws.Range("J91").Value = 25000
$ws.Range("L91").Value = 25000
$ws.Range("N91").Value = -27808
$ws.Range("H113").Value = 407.33334
$ws.Range("I113").Value = 118.666664
$ws.Range("K113").Value = 355.999992
$ws.Range("M113").Value = 1814.000008
$ws.Range("H114").Value = 57500
$ws.Range("J114").Value = 57500
$ws.Range("L114").Value = 57500
$ws.Range("N114").Value = -66178
$ws.Range("H124").Value = 23147
$ws.Range("J124").Value = 23147
$ws.Range("L124").Value = 23147
$ws.Range("N124").Value = -32967
$ws.Range("H132").Value = 2326.8
$ws.Range("I132").Value = 2127.318
$ws.Range("K132").Value = 6381.954000000001
$ws.Range("M132").Value = -3851.954000000001
$ws.Range("H135").Value = 99999
$ws.Range("J135").Value = 99999
$ws.Range("L135").Value = 99999
$ws.Range("N135").Value = -110139
